# Scheduled runner update: refresh recorded market-board prices / leve
# profit calculations (columns H:N) across several crafting-job sheets.
# Column layout per row: H=currentAveragePrice, I=currentAveragePriceNQ,
# J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ,
# N=LeveProfitHQ.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value = 3036.5264
$ws.Cells.Item(58, 9).Value = 215
$ws.Cells.Item(58, 10).Value = 3193.2778
$ws.Cells.Item(58, 11).Value = 645
$ws.Cells.Item(58, 12).Value = 9579.8334
$ws.Cells.Item(58, 13).Value = -495
$ws.Cells.Item(58, 14).Value = -9879.8334
$ws.Cells.Item(62, 8).Value = 2023.4231
$ws.Cells.Item(62, 9).Value = 2300.1538
$ws.Cells.Item(62, 11).Value = 2300.1538
$ws.Cells.Item(62, 13).Value = -1676.1538
$ws.Cells.Item(65, 8).Value = 2023.4231
$ws.Cells.Item(65, 9).Value = 2300.1538
$ws.Cells.Item(65, 11).Value = 11500.769
$ws.Cells.Item(65, 13).Value = -8380.769
$ws.Cells.Item(107, 8).Value = 14706389
$ws.Cells.Item(107, 9).Value = 16667107
$ws.Cells.Item(107, 11).Value = 16667107
$ws.Cells.Item(107, 13).Value = -16665187
$ws.Cells.Item(112, 8).Value = 21164982
$ws.Cells.Item(112, 10).Value = 22858148
$ws.Cells.Item(112, 12).Value = 68574444
$ws.Cells.Item(112, 14).Value = -68576660
$ws.Cells.Item(131, 8).Value = 2751.5
$ws.Cells.Item(131, 9).Value = 485.83334
$ws.Cells.Item(131, 11).Value = 1457.50002
$ws.Cells.Item(131, 13).Value = 3582.49998
$ws.Cells.Item(135, 8).Value = 1924.644
$ws.Cells.Item(135, 9).Value = 1699.7046
$ws.Cells.Item(135, 10).Value = 2584.4666
$ws.Cells.Item(135, 11).Value = 15297.3414
$ws.Cells.Item(135, 12).Value = 23260.1994
$ws.Cells.Item(135, 13).Value = -12762.3414
$ws.Cells.Item(135, 14).Value = -28330.1994
$ws.Cells.Item(137, 8).Value = 1423.9298
$ws.Cells.Item(137, 9).Value = 1112.2094
$ws.Cells.Item(137, 10).Value = 2381.3572
$ws.Cells.Item(137, 11).Value = 3336.6282
$ws.Cells.Item(137, 12).Value = 7144.071599999999
$ws.Cells.Item(137, 13).Value = -786.6282000000001
$ws.Cells.Item(137, 14).Value = -12244.0716

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 540.1539
$ws.Cells.Item(2, 9).Value = 426.1111
$ws.Cells.Item(2, 10).Value = 1908.6666
$ws.Cells.Item(2, 11).Value = 426.1111
$ws.Cells.Item(2, 12).Value = 1908.6666
$ws.Cells.Item(2, 13).Value = -313.1111
$ws.Cells.Item(2, 14).Value = -2134.6666
$ws.Cells.Item(26, 8).Value = 4235.2856
$ws.Cells.Item(26, 9).Value = 1661.75
$ws.Cells.Item(26, 10).Value = 7666.6665
$ws.Cells.Item(26, 11).Value = 1661.75
$ws.Cells.Item(26, 12).Value = 7666.6665
$ws.Cells.Item(26, 13).Value = -1331.75
$ws.Cells.Item(26, 14).Value = -8326.666499999999
$ws.Cells.Item(32, 8).Value = 9756.806
$ws.Cells.Item(32, 9).Value = 7411.5635
$ws.Cells.Item(32, 11).Value = 7411.5635
$ws.Cells.Item(32, 13).Value = -7124.5635
$ws.Cells.Item(33, 8).Value = 5750
$ws.Cells.Item(33, 10).Value = 5750
$ws.Cells.Item(33, 12).Value = 5750
$ws.Cells.Item(33, 14).Value = -6408
$ws.Cells.Item(116, 8).Value = 540.1539
$ws.Cells.Item(116, 9).Value = 426.1111
$ws.Cells.Item(116, 10).Value = 1908.6666
$ws.Cells.Item(116, 11).Value = 426.1111
$ws.Cells.Item(116, 12).Value = 1908.6666
$ws.Cells.Item(116, 13).Value = 1867.8889
$ws.Cells.Item(116, 14).Value = -6496.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 540.1539
$ws.Cells.Item(3, 9).Value = 426.1111
$ws.Cells.Item(3, 10).Value = 1908.6666
$ws.Cells.Item(3, 11).Value = 426.1111
$ws.Cells.Item(3, 12).Value = 1908.6666
$ws.Cells.Item(3, 13).Value = -312.1111
$ws.Cells.Item(3, 14).Value = -2136.6666
$ws.Cells.Item(20, 8).Value = 12779.182
$ws.Cells.Item(20, 10).Value = 50981.8
$ws.Cells.Item(20, 12).Value = 50981.8
$ws.Cells.Item(20, 14).Value = -51475.8
$ws.Cells.Item(105, 8).Value = 14138.474
$ws.Cells.Item(105, 9).Value = 31572.715
$ws.Cells.Item(105, 11).Value = 31572.715
$ws.Cells.Item(105, 13).Value = -29825.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1815.2
$ws.Cells.Item(16, 9).Value = 1479
$ws.Cells.Item(16, 11).Value = 1479
$ws.Cells.Item(16, 13).Value = -1192
$ws.Cells.Item(31, 8).Value = 2811.8452
$ws.Cells.Item(31, 9).Value = 1748.9302
$ws.Cells.Item(31, 10).Value = 3926.6099
$ws.Cells.Item(31, 11).Value = 1748.9302
$ws.Cells.Item(31, 12).Value = 3926.6099
$ws.Cells.Item(31, 13).Value = -1453.9302
$ws.Cells.Item(31, 14).Value = -4516.609899999999
$ws.Cells.Item(34, 8).Value = 2811.8452
$ws.Cells.Item(34, 9).Value = 1748.9302
$ws.Cells.Item(34, 10).Value = 3926.6099
$ws.Cells.Item(34, 11).Value = 1748.9302
$ws.Cells.Item(34, 12).Value = 3926.6099
$ws.Cells.Item(34, 13).Value = -1546.9302
$ws.Cells.Item(34, 14).Value = -4330.609899999999
$ws.Cells.Item(35, 8).Value = 4380.3335
$ws.Cells.Item(35, 9).Value = 3008.3333
$ws.Cells.Item(35, 10).Value = 5752.3335
$ws.Cells.Item(35, 11).Value = 3008.3333
$ws.Cells.Item(35, 12).Value = 5752.3335
$ws.Cells.Item(35, 13).Value = -2714.3333
$ws.Cells.Item(35, 14).Value = -6340.3335
$ws.Cells.Item(105, 8).Value = 2624.68
$ws.Cells.Item(105, 9).Value = 2573.625
$ws.Cells.Item(105, 11).Value = 2573.625
$ws.Cells.Item(105, 13).Value = -826.625
$ws.Cells.Item(113, 8).Value = 1815.2
$ws.Cells.Item(113, 9).Value = 1479
$ws.Cells.Item(113, 11).Value = 1479
$ws.Cells.Item(113, 13).Value = 691

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1062
$ws.Cells.Item(34, 9).Value = 928.3333
$ws.Cells.Item(34, 10).Value = 1151.1111
$ws.Cells.Item(34, 11).Value = 2784.9999
$ws.Cells.Item(34, 12).Value = 3453.3333
$ws.Cells.Item(34, 13).Value = -2700.9999
$ws.Cells.Item(34, 14).Value = -3621.3333
$ws.Cells.Item(39, 8).Value = 1349.5
$ws.Cells.Item(39, 10).Value = 1349.5
$ws.Cells.Item(39, 12).Value = 4048.5
$ws.Cells.Item(39, 14).Value = -4636.5
$ws.Cells.Item(55, 8).Value = 2751.8518
$ws.Cells.Item(55, 9).Value = 500
$ws.Cells.Item(55, 10).Value = 2838.4614
$ws.Cells.Item(55, 11).Value = 1500
$ws.Cells.Item(55, 12).Value = 8515.3842
$ws.Cells.Item(55, 13).Value = -1323
$ws.Cells.Item(55, 14).Value = -8869.3842
$ws.Cells.Item(109, 8).Value = 1868.75
$ws.Cells.Item(109, 9).Value = 350
$ws.Cells.Item(109, 11).Value = 1050
$ws.Cells.Item(109, 13).Value = -10
$ws.Cells.Item(129, 8).Value = 1031.75
$ws.Cells.Item(129, 9).Value = 500
$ws.Cells.Item(129, 10).Value = 1125.5883
$ws.Cells.Item(129, 11).Value = 1500
$ws.Cells.Item(129, 12).Value = 3376.7649
$ws.Cells.Item(129, 13).Value = 3500
$ws.Cells.Item(129, 14).Value = -13376.7649
$ws.Cells.Item(132, 8).Value = 1102.6666
$ws.Cells.Item(132, 9).Value = 582.8333
$ws.Cells.Item(132, 10).Value = 1622.5
$ws.Cells.Item(132, 11).Value = 5245.4997
$ws.Cells.Item(132, 12).Value = 14602.5
$ws.Cells.Item(132, 13).Value = -2715.4997
$ws.Cells.Item(132, 14).Value = -19662.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1379.4762
$ws.Cells.Item(97, 9).Value = 1379.4762
$ws.Cells.Item(97, 11).Value = 1379.4762
$ws.Cells.Item(97, 13).Value = -883.4762000000001
$ws.Cells.Item(123, 8).Value = 10190.125
$ws.Cells.Item(123, 9).Value = 9000
$ws.Cells.Item(123, 10).Value = 10241.869
$ws.Cells.Item(123, 11).Value = 9000
$ws.Cells.Item(123, 12).Value = 10241.869
$ws.Cells.Item(123, 13).Value = -6550
$ws.Cells.Item(123, 14).Value = -15141.869

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1498.4615
$ws.Cells.Item(16, 9).Value = 1539.8334
$ws.Cells.Item(16, 11).Value = 1539.8334
$ws.Cells.Item(16, 13).Value = -1369.8334
$ws.Cells.Item(61, 8).Value = 1766.1578
$ws.Cells.Item(61, 9).Value = 1459.7273
$ws.Cells.Item(61, 10).Value = 2187.5
$ws.Cells.Item(61, 11).Value = 1459.7273
$ws.Cells.Item(61, 12).Value = 2187.5
$ws.Cells.Item(61, 13).Value = -1257.7273
$ws.Cells.Item(61, 14).Value = -2591.5
$ws.Cells.Item(100, 8).Value = 1087.84
$ws.Cells.Item(100, 9).Value = 995.05554
$ws.Cells.Item(100, 10).Value = 1326.4286
$ws.Cells.Item(100, 11).Value = 995.05554
$ws.Cells.Item(100, 12).Value = 1326.4286
$ws.Cells.Item(100, 13).Value = -454.05554
$ws.Cells.Item(100, 14).Value = -2408.4286
$ws.Cells.Item(113, 8).Value = 1766.1578
$ws.Cells.Item(113, 9).Value = 1459.7273
$ws.Cells.Item(113, 10).Value = 2187.5
$ws.Cells.Item(113, 11).Value = 1459.7273
$ws.Cells.Item(113, 12).Value = 2187.5
$ws.Cells.Item(113, 13).Value = 710.2727
$ws.Cells.Item(113, 14).Value = -6527.5
$ws.Cells.Item(133, 8).Value = 130108.664
$ws.Cells.Item(133, 10).Value = 130108.664
$ws.Cells.Item(133, 12).Value = 130108.664
$ws.Cells.Item(133, 14).Value = -135168.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 100000500
$ws.Cells.Item(107, 9).Value = 111111520
$ws.Cells.Item(107, 11).Value = 333334560
$ws.Cells.Item(107, 13).Value = -333332640
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 1105.7391
$ws.Cells.Item(113, 9).Value = 840.7941
$ws.Cells.Item(113, 10).Value = 1856.4166
$ws.Cells.Item(113, 11).Value = 2522.3823
$ws.Cells.Item(113, 12).Value = 5569.2498
$ws.Cells.Item(113, 13).Value = -352.3822999999998
$ws.Cells.Item(113, 14).Value = -9909.2498
